# Update attendance summary columns (D-H) for rows 3-18 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Total Attendance Count (D), Real (E), Duplicate (F), Invalid (G), Absent (H)
$values = @{
    3  = @(0, 0, 0, 1, 1)
    4  = @(2, 1, 1, 0, 0)
    5  = @(0, 0, 0, 0, 1)
    6  = @(1, 1, 0, 0, 0)
    7  = @(0, 0, 0, 0, 1)
    8  = @(0, 0, 0, 0, 1)
    9  = @(0, 0, 0, 0, 1)
    10 = @(0, 0, 0, 0, 1)
    11 = @(0, 0, 0, 0, 1)
    12 = @(1, 1, 0, 0, 0)
    13 = @(1, 1, 0, 0, 0)
    14 = @(0, 0, 0, 0, 1)
    15 = @(0, 0, 0, 0, 1)
    16 = @(0, 0, 0, 0, 1)
    17 = @(0, 0, 0, 0, 1)
    18 = @(0, 0, 0, 0, 1)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("D$row").Value = $rowValues[0]
    $ws.Range("E$row").Value = $rowValues[1]
    $ws.Range("F$row").Value = $rowValues[2]
    $ws.Range("G$row").Value = $rowValues[3]
    $ws.Range("H$row").Value = $rowValues[4]
}
